$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = "62.663.92"
$ws.Range("E2").Value = "  +4.03%  "
$ws.Range("D3").Value = "3.344.87"
$ws.Range("E3").Value = "  +4.16%  "
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.06%  "
Set-TextValue "D5" "558.80"
$ws.Range("E5").Value = "  +3.74%  "
Set-TextValue "D6" "152.38"
$ws.Range("E6").Value = "  +4.61%  "
Set-TextValue "D7" "0.999"
$ws.Range("E7").Value = "  +0.00%  "
Set-TextValue "D8" "0.533"
$ws.Range("E8").Value = "  -0.72%  "
$ws.Range("E9").Value = "  +1.43%  "
$ws.Range("E10").Value = "  +3.77%  "
Set-TextValue "D11" "0.437"
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("D12").Value = "3.920.09"
$ws.Range("E12").Value = "  +4.18%  "
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D14" "26.94"
$ws.Range("E14").Value = "  +2.82%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D15" "0.0000180"
$ws.Range("E15").Value = "  +2.99%  "
$ws.Range("D16").Value = "62.654.76"
$ws.Range("E16").Value = "  +3.93%  "
$ws.Range("D17").Value = "3.318.36"
$ws.Range("E17").Value = "  +3.26%  "
Set-TextValue "D18" "6.42"
Set-TextValue "D19" "13.83"
$ws.Range("E19").Value = "  +4.86%  "
Set-TextValue "D20" "386.51"
$ws.Range("E20").Value = "  +1.27%  "
Set-TextValue "D21" "8.31"
$ws.Range("E21").Value = "  -0.60%  "
Set-TextValue "D22" "0.999"
$ws.Range("E22").Value = "  -0.10%  "
Set-TextValue "D23" "0.538"
$ws.Range("E23").Value = "  +1.66%  "
Set-TextValue "D24" "69.99"
$ws.Range("E24").Value = "  -0.41%  "
$ws.Range("E25").Value = "  +5.20%  "
Set-TextValue "D26" "8.83"
$ws.Range("E26").Value = "  -0.97%  "
$ws.Range("D27").Value = "0.0₃0952"
$ws.Range("E27").Value = "  +4.54%  "
$ws.Range("E28").Value = "  +0.11%  "
Set-TextValue "D29" "6.57"
$ws.Range("E29").Value = "  +6.36%  "
Set-TextValue "D30" "1.98"
$ws.Range("E30").Value = "  +3.35%  "
$ws.Range("E31").Value = "  +2.78%  "
Set-TextValue "D32" "22.94"
$ws.Range("E32").Value = "  +2.26%  "
$ws.Range("E33").Value = "  +6.25%  "
Set-TextValue "D34" "6.68"
$ws.Range("E34").Value = "  +0.81%  "
Set-TextValue "D35" "160.33"
$ws.Range("E35").Value = "  +2.06%  "
$ws.Range("E36").Value = "  +7.36%  "
$ws.Range("E37").Value = "  +11.91%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D38" "0.0747"
$ws.Range("E38").Value = "  +5.45%  "
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D39" "26.70"
$ws.Range("E39").Value = "  +3.03%  "
$ws.Range("D40").Value = "2.830.79"
$ws.Range("E40").Value = "  +1.23%  "
Set-TextValue "D41" "0.0312"
$ws.Range("E41").Value = "  +8.51%  "
$ws.Range("E42").Value = "  +3.42%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D43" "40.56"
$ws.Range("E43").Value = "  +1.09%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D44" "4.27"
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("B45").Value = "ONDO"
$ws.Range("C45").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue "D45" "1.03"
$ws.Range("E45").Value = "  +2.89%  "
$ws.Range("B46").Value = "RenzoRestakedETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D46").Value = "3.388.09"
$ws.Range("E46").Value = "  +4.22%  "
Set-TextValue "D47" "21.96"
$ws.Range("E47").Value = "  +5.27%  "
Set-TextValue "D48" "0.104"
$ws.Range("E48").Value = "  +0.25%  "
Set-TextValue "D49" "6.28"
$ws.Range("E49").Value = "  +1.56%  "
Set-TextValue "D50" "0.797"
$ws.Range("E50").Value = "  -1.39%  "
Set-TextValue "D51" "283.08"
$ws.Range("E51").Value = "  +3.90%  "
